# Fix [ToC] in Word master and html
#
# The placeholder paragraph "[ToC]" was split across three runs
# (with spell-check proofErr markers bracketing the "ToC" run). Replace
# it with a single, unbroken run reading "[ToC]" and make the paragraph
# bold (bold applied at the paragraph-mark level, i.e. w:pPr/w:rPr).

$d = $word.ActiveDocument

# Locate the "[ToC]" paragraph robustly (by content, not a hard-coded index).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.Trim() -eq "[ToC]") {
        $target = $candidate
        break
    }
}

if ($target -ne $null) {
    # Pull the paragraph's own WordOpenXML so we can keep its existing
    # w14:paraId / rsid* attributes instead of clobbering them.
    $openXml = $target.Range.WordOpenXML
    $openTag = "<w:p>"
    if ($openXml -match "(<w:p [^>]*>)") {
        $openTag = $matches[1]
    }
    $attrs = ""
    if ($openTag -match "^<w:p\s+(.*)>$") {
        $attrs = " " + $matches[1]
    }

    $newParaXml = "<w:p xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main""" +
        " xmlns:w14=""http://schemas.microsoft.com/office/word/2010/wordml""" + $attrs + ">" +
        "<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>" +
        "<w:r><w:t>[ToC]</w:t></w:r>" +
        "</w:p>"

    $target.Range.InsertXML($newParaXml)
}
